$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8663.16
$ws.Range("J19").Value = 10157.523
$ws.Range("L19").Value = 10157.523
$ws.Range("N19").Value = -10507.523

$ws.Range("H21").Value = 18913.715
$ws.Range("J21").Value = 18913.715
$ws.Range("L21").Value = 18913.715
$ws.Range("N21").Value = -19849.715

$ws.Range("H23").Value = 18913.715
$ws.Range("J23").Value = 18913.715
$ws.Range("L23").Value = 18913.715
$ws.Range("N23").Value = -19381.715

$ws.Range("H32").Value = 522.1539
$ws.Range("I32").Value = 422.75
$ws.Range("J32").Value = 566.3333
$ws.Range("K32").Value = 422.75
$ws.Range("L32").Value = 566.3333
$ws.Range("M32").Value = -96.75
$ws.Range("N32").Value = -1218.3333

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H51").Value = 1433.3334
$ws.Range("I51").Value = 1400
$ws.Range("J51").Value = 1500
$ws.Range("K51").Value = 1400
$ws.Range("L51").Value = 1500
$ws.Range("M51").Value = -916
$ws.Range("N51").Value = -2468

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H92").Value = 23810222
$ws.Range("I92").Value = 25641470
$ws.Range("K92").Value = 25641470
$ws.Range("M92").Value = -25640222

$ws.Range("H94").Value = 676.25
$ws.Range("I94").Value = 676.25
$ws.Range("K94").Value = 676.25
$ws.Range("M94").Value = -225.25

$ws.Range("H107").Value = 1062.0416
$ws.Range("I107").Value = 1247.5264
$ws.Range("J107").Value = 357.2
$ws.Range("K107").Value = 1247.5264
$ws.Range("L107").Value = 357.2
$ws.Range("M107").Value = 672.4736
$ws.Range("N107").Value = -4197.2

$ws.Range("H113").Value = 2186.7144
$ws.Range("I113").Value = 1770.3846
$ws.Range("J113").Value = 2863.25
$ws.Range("K113").Value = 1770.3846
$ws.Range("L113").Value = 2863.25
$ws.Range("M113").Value = 1483.6154
$ws.Range("N113").Value = -9371.25

$ws.Range("H116").Value = 2629.8235
$ws.Range("I116").Value = 2556.111
$ws.Range("J116").Value = 2712.75
$ws.Range("K116").Value = 2556.111
$ws.Range("L116").Value = 2712.75
$ws.Range("M116").Value = 885.8890000000001
$ws.Range("N116").Value = -9596.75

$ws.Range("H132").Value = 3070.6177
$ws.Range("I132").Value = 2860.6667
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 8582.000100000001
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -6052.000100000001
$ws.Range("N132").Value = -35057

$ws.Range("H133").Value = 42400
$ws.Range("J133").Value = 42400
$ws.Range("L133").Value = 42400
$ws.Range("N133").Value = -52520

$ws.Range("H134").Value = 62156
$ws.Range("J134").Value = 62156
$ws.Range("L134").Value = 62156
$ws.Range("N134").Value = -72296

$ws.Range("H136").Value = 39212.727
$ws.Range("J136").Value = 39212.727
$ws.Range("L136").Value = 39212.727
$ws.Range("N136").Value = -49412.727

$ws.Range("H137").Value = 2008.6154
$ws.Range("I137").Value = 2231.7
$ws.Range("J137").Value = 1869.1875
$ws.Range("K137").Value = 6695.099999999999
$ws.Range("L137").Value = 5607.5625
$ws.Range("M137").Value = -4145.099999999999
$ws.Range("N137").Value = -10707.5625

$ws.Range("H138").Value = 2137.45
$ws.Range("I138").Value = 1081.2307
$ws.Range("J138").Value = 2295.276
$ws.Range("K138").Value = 3243.6921
$ws.Range("L138").Value = 6885.828
$ws.Range("M138").Value = 1896.3079
$ws.Range("N138").Value = -17165.828

$ws.Range("H139").Value = 50780
$ws.Range("J139").Value = 50780
$ws.Range("L139").Value = 50780
$ws.Range("N139").Value = -61060

$ws.Range("H140").Value = 50780
$ws.Range("J140").Value = 50780
$ws.Range("L140").Value = 50780
$ws.Range("N140").Value = -61140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1021215.4
$ws.Range("I32").Value = 1219683.8
$ws.Range("J32").Value = 28873.092
$ws.Range("K32").Value = 1219683.8
$ws.Range("L32").Value = 28873.092
$ws.Range("M32").Value = -1219396.8
$ws.Range("N32").Value = -29447.092

$ws.Range("H37").Value = 9994.333000000001
$ws.Range("I37").Value = 2900
$ws.Range("J37").Value = 12021.286
$ws.Range("K37").Value = 2900
$ws.Range("L37").Value = 12021.286
$ws.Range("M37").Value = -2627
$ws.Range("N37").Value = -12567.286

$ws.Range("H44").Value = 29999
$ws.Range("J44").Value = 29999
$ws.Range("L44").Value = 29999
$ws.Range("N44").Value = -30975

$ws.Range("I74").Value = 1173.5
$ws.Range("K74").Value = 1173.5
$ws.Range("M74").Value = -299.5

$ws.Range("I77").Value = 1173.5
$ws.Range("K77").Value = 5867.5
$ws.Range("M77").Value = -1499.5

$ws.Range("H122").Value = 168785.67
$ws.Range("I122").Value = 334000
$ws.Range("J122").Value = 3571.3333
$ws.Range("K122").Value = 1002000
$ws.Range("L122").Value = 10713.9999
$ws.Range("M122").Value = -999550
$ws.Range("N122").Value = -15613.9999

$ws.Range("H132").Value = 1640414.6
$ws.Range("I132").Value = 3035
$ws.Range("J132").Value = 3209570.2
$ws.Range("K132").Value = 9105
$ws.Range("L132").Value = 9628710.600000001
$ws.Range("M132").Value = -6575
$ws.Range("N132").Value = -9633770.600000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 19999
$ws.Range("J35").Value = 19999
$ws.Range("L35").Value = 19999
$ws.Range("N35").Value = -20619

$ws.Range("H134").Value = 2461.4736
$ws.Range("I134").Value = 2864.7273
$ws.Range("K134").Value = 8594.1819
$ws.Range("M134").Value = -6059.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5855.4604
$ws.Range("I31").Value = 1774.9565
$ws.Range("J31").Value = 8201.75
$ws.Range("K31").Value = 1774.9565
$ws.Range("L31").Value = 8201.75
$ws.Range("M31").Value = -1479.9565
$ws.Range("N31").Value = -8791.75

$ws.Range("H34").Value = 5855.4604
$ws.Range("I34").Value = 1774.9565
$ws.Range("J34").Value = 8201.75
$ws.Range("K34").Value = 1774.9565
$ws.Range("L34").Value = 8201.75
$ws.Range("M34").Value = -1572.9565
$ws.Range("N34").Value = -8605.75

$ws.Range("H99").Value = 1840.2273
$ws.Range("I99").Value = 1340.4286
$ws.Range("J99").Value = 2073.4666
$ws.Range("K99").Value = 1340.4286
$ws.Range("L99").Value = 2073.4666
$ws.Range("M99").Value = 157.5714
$ws.Range("N99").Value = -5069.4666

$ws.Range("H126").Value = 1840.2273
$ws.Range("I126").Value = 1340.4286
$ws.Range("J126").Value = 2073.4666
$ws.Range("K126").Value = 4021.2858
$ws.Range("L126").Value = 6220.399800000001
$ws.Range("M126").Value = -1551.2858
$ws.Range("N126").Value = -11160.3998

$ws.Range("H132").Value = 1631.2858
$ws.Range("I132").Value = 1090.2632
$ws.Range("J132").Value = 2773.4443
$ws.Range("K132").Value = 3270.7896
$ws.Range("L132").Value = 8320.332900000001
$ws.Range("M132").Value = -740.7896000000001
$ws.Range("N132").Value = -13380.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 4265.1875
$ws.Range("I49").Value = 3937.875
$ws.Range("J49").Value = 4592.5
$ws.Range("K49").Value = 11813.625
$ws.Range("L49").Value = 13777.5
$ws.Range("M49").Value = -11657.625
$ws.Range("N49").Value = -14089.5

$ws.Range("H56").Value = 5340
$ws.Range("I56").Value = 5340
$ws.Range("K56").Value = 5340
$ws.Range("M56").Value = -4810

$ws.Range("H131").Value = 3412.434
$ws.Range("I131").Value = 626
$ws.Range("J131").Value = 3702.6875
$ws.Range("K131").Value = 1878
$ws.Range("L131").Value = 11108.0625
$ws.Range("M131").Value = 3162
$ws.Range("N131").Value = -21188.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 26812592
$ws.Range("I80").Value = 42418500
$ws.Range("J80").Value = 59606.145
$ws.Range("K80").Value = 42418500
$ws.Range("L80").Value = 59606.145
$ws.Range("M80").Value = -42417502
$ws.Range("N80").Value = -61602.145

$ws.Range("H83").Value = 26812592
$ws.Range("I83").Value = 42418500
$ws.Range("J83").Value = 59606.145
$ws.Range("K83").Value = 212092500
$ws.Range("L83").Value = 298030.725
$ws.Range("M83").Value = -212087508
$ws.Range("N83").Value = -308014.725

$ws.Range("H102").Value = 2352.2222
$ws.Range("I102").Value = 2344
$ws.Range("K102").Value = 2344
$ws.Range("M102").Value = -722

$ws.Range("H132").Value = 3081.1025
$ws.Range("I132").Value = 2698.4517
$ws.Range("J132").Value = 4563.875
$ws.Range("K132").Value = 8095.355100000001
$ws.Range("L132").Value = 13691.625
$ws.Range("M132").Value = -5565.355100000001
$ws.Range("N132").Value = -18751.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3999.3044
$ws.Range("I7").Value = 4433.143
$ws.Range("J7").Value = 3324.4443
$ws.Range("K7").Value = 4433.143
$ws.Range("L7").Value = 3324.4443
$ws.Range("M7").Value = -4321.143
$ws.Range("N7").Value = -3548.4443

$ws.Range("H31").Value = 2673.2856
$ws.Range("I31").Value = 832.2857
$ws.Range("J31").Value = 4514.2856
$ws.Range("K31").Value = 832.2857
$ws.Range("L31").Value = 4514.2856
$ws.Range("M31").Value = -584.2857
$ws.Range("N31").Value = -5010.2856

$ws.Range("H40").Value = 335268
$ws.Range("I40").Value = 335268
$ws.Range("K40").Value = 335268
$ws.Range("M40").Value = -335132

$ws.Range("H100").Value = 4965
$ws.Range("I100").Value = 4953.3335
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 4953.3335
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = -4412.3335
$ws.Range("N100").Value = -6082

$ws.Range("H126").Value = 3999.3044
$ws.Range("I126").Value = 4433.143
$ws.Range("J126").Value = 3324.4443
$ws.Range("K126").Value = 13299.429
$ws.Range("L126").Value = 9973.332900000001
$ws.Range("M126").Value = -10829.429
$ws.Range("N126").Value = -14913.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 54041.668
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 64250
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 64250
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -65498

$ws.Range("H65").Value = 54041.668
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 64250
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 321250
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -327490

$ws.Range("H132").Value = 9117271
$ws.Range("I132").Value = 3923.1428
$ws.Range("J132").Value = 16205431
$ws.Range("K132").Value = 11769.4284
$ws.Range("L132").Value = 48616293
$ws.Range("M132").Value = -9239.428400000001
$ws.Range("N132").Value = -48621353
